$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.973.07'
$ws.Range("E2").Value = '  +2.44%  '

$ws.Range("D3").Value = '3.739.66'
$ws.Range("E3").Value = '  -0.90%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = "'601.29"
$ws.Range("E5").Value = '  +1.83%  '

$ws.Range("D6").Value = "'168.59"
$ws.Range("E6").Value = '  -1.32%  '

$ws.Range("D7").Value = '3.736.78'
$ws.Range("E7").Value = '  -0.93%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").Value = "'0.532"
$ws.Range("E9").Value = '  +2.47%  '

$ws.Range("D10").Value = "'0.165"
$ws.Range("E10").Value = '  +4.32%  '

$ws.Range("E11").Value = '  +1.21%  '

$ws.Range("D13").Value = "'38.26"
$ws.Range("E13").Value = '  +1.96%  '

$ws.Range("E14").Value = '  +0.57%  '

$ws.Range("D15").Value = '4.359.62'
$ws.Range("E15").Value = '  -0.97%  '

$ws.Range("D16").Value = '3.729.44'
$ws.Range("E16").Value = '  -1.06%  '

$ws.Range("D17").Value = '68.925.98'
$ws.Range("E17").Value = '  +2.14%  '

$ws.Range("D18").Value = "'7.26"

$ws.Range("E19").Value = '  +0.36%  '

$ws.Range("D20").Value = "'17.21"
$ws.Range("E20").Value = '  +7.95%  '

$ws.Range("D21").Value = "'497.66"
$ws.Range("E21").Value = '  +2.54%  '

$ws.Range("D22").Value = "'9.53"
$ws.Range("E22").Value = '  +4.37%  '

$ws.Range("D23").Value = "'0.723"
$ws.Range("E23").Value = '  +0.72%  '

$ws.Range("D24").Value = "'84.82"
$ws.Range("E24").Value = '  +1.25%  '

$ws.Range("E25").Value = '  -2.00%  '

$ws.Range("E26").Value = '  +1.71%  '

$ws.Range("E27").Value = '  +1.37%  '

$ws.Range("D28").Value = "'10.11"
$ws.Range("E28").Value = '  -0.18%  '

$ws.Range("E29").Value = '  +0.03%  '

$ws.Range("D30").Value = "'2.94"
$ws.Range("E30").Value = '  +1.54%  '

$ws.Range("D31").Value = "'2.42"
$ws.Range("E31").Value = '  +1.62%  '

$ws.Range("D32").Value = "'7.93"
$ws.Range("E32").Value = '  +2.25%  '

$ws.Range("D33").Value = "'31.74"
$ws.Range("E33").Value = '  -1.66%  '

$ws.Range("D34").Value = '3.876.03'
$ws.Range("E34").Value = '  -0.96%  '

$ws.Range("E35").Value = '  +1.00%  '

$ws.Range("D36").Value = '3.670.55'
$ws.Range("E36").Value = '  -1.13%  '

$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = '  -0.01%  '

$ws.Range("E38").Value = '  +1.04%  '

$ws.Range("E39").Value = '  +1.28%  '

$ws.Range("D40").Value = "'0.134"
$ws.Range("E40").Value = '  -1.09%  '

$ws.Range("D41").Value = "'0.324"
$ws.Range("E41").Value = '  +0.85%  '

$ws.Range("D42").Value = "'437.37"
$ws.Range("E42").Value = '  -2.48%  '

$ws.Range("D43").Value = "'49.00"
$ws.Range("E43").Value = '  +0.52%  '

$ws.Range("E44").Value = '  -0.21%  '

$ws.Range("D45").Value = "'2.88"
$ws.Range("E45").Value = '  +1.68%  '

$ws.Range("D46").Value = "'8.39"
$ws.Range("E46").Value = '  +2.02%  '

$ws.Range("E47").Value = '  +0.00%  '

$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").Value = "'40.51"
$ws.Range("E48").Value = '  -2.08%  '

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = "'143.68"
$ws.Range("E49").Value = '  +2.30%  '

$ws.Range("D50").Value = "'0.0352"
$ws.Range("E50").Value = '  +1.43%  '

$ws.Range("D51").Value = '2.746.66'
$ws.Range("E51").Value = '  -2.65%  '
